$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells stay plain text (not auto-converted to numbers)
# by forcing a text number format before assigning values, then resetting the
# cell style back to Normal so no stray formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.763.96'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.850.31'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '315.29'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').Value = '0.4273'
$ws.Range('E7').Value = '  -2.10%  '
$ws.Range('D8').Value = '0.3668'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').Value = '45.15'
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('D10').Value = '0.07330'
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('D11').Value = '0.8936'
$ws.Range('E11').Value = '  -4.44%  '
$ws.Range('D12').Value = '20.91'
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('D13').Value = '1.834.77'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = '6.586'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = '5.358'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = '0.06920'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -2.90%  '
$ws.Range('D19').Value = '0.000008925'
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').Value = '15.52'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '27.758.30'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('D23').Value = '4.995'
$ws.Range('E23').Value = '  -2.55%  '
$ws.Range('D24').Value = '10.64'
$ws.Range('E24').Value = '  -3.38%  '
$ws.Range('D25').Value = '2.089.68'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('D26').Value = '1.958'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = '153.63'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').Value = '18.90'
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('D29').Value = '121.03'
$ws.Range('E29').Value = '  +6.96%  '
$ws.Range('D30').Value = '5.269'
$ws.Range('E30').Value = '  -3.89%  '
$ws.Range('D31').Value = '1.917'
$ws.Range('E31').Value = '  +11.64%  '
$ws.Range('D32').Value = '0.08945'
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('D33').Value = '0.7729'
$ws.Range('E33').Value = '  -6.00%  '
$ws.Range('D34').Value = '4.602'
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('D35').Value = '2.980'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '1.107'
$ws.Range('E36').Value = '  -6.01%  '
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = '1.100'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').Value = '0.05396'
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('D40').Value = '0.01960'
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').Value = '2.821'
$ws.Range('E41').Value = '  -4.84%  '
$ws.Range('D42').Value = '6.921'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').Value = '0.5127'
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').Value = '0.1665'
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').Value = '8.296'
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('D46').Value = '0.06597'
$ws.Range('E46').Value = '  -2.48%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.4781'
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '10.48'
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('D49').Value = '104.78'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '1.639'
$ws.Range('E51').Value = '  -2.41%  '

$priceRange.Style = "Normal"
